$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.150083774525282
$ws.Range("D2").Value = 3.237573490718583
$ws.Range("E2").Value = 10.65779238487089
$ws.Range("F2").Value = 27.42418067053856
$ws.Range("G2").Value = 38.6082390399357
$ws.Range("H2").Value = 13.93665405262289
$ws.Range("I2").Value = 17.39779600370851
$ws.Range("M2").Value = 19.01097732685402
$ws.Range("N2").Value = 16.74216291963044

$ws.Range("B3").Value = 6.96126610514583
$ws.Range("D3").Value = 3.221632694131318
$ws.Range("E3").Value = 10.75199199973205
$ws.Range("F3").Value = 26.45221855932836
$ws.Range("G3").Value = 36.65244489645163
$ws.Range("H3").Value = 13.71826460683177
$ws.Range("I3").Value = 17.46551527722893
$ws.Range("M3").Value = 18.20532322585688
$ws.Range("N3").Value = 16.71726560329046

$ws.Range("B4").Value = 6.843164313807258
$ws.Range("D4").Value = 3.213198901088585
$ws.Range("E4").Value = 10.81259673085782
$ws.Range("F4").Value = 25.84858229052279
$ws.Range("G4").Value = 35.41245292027211
$ws.Range("H4").Value = 13.5877488252076
$ws.Range("I4").Value = 17.51090699389515
$ws.Range("M4").Value = 17.69575121732931
$ws.Range("N4").Value = 16.70409753502454

$ws.Range("B5").Value = 6.794558456863385
$ws.Range("D5").Value = 3.210102356666742
$ws.Range("E5").Value = 10.83799223349365
$ws.Range("F5").Value = 25.60132012339781
$ws.Range("G5").Value = 34.89807123698314
$ws.Range("H5").Value = 13.53552980815524
$ws.Range("I5").Value = 17.53036315264098
$ws.Range("M5").Value = 17.48465987682765
$ws.Range("N5").Value = 16.69926957048985

$ws.Range("B6").Value = 6.786460617431629
$ws.Range("D6").Value = 3.209608713459607
$ws.Range("E6").Value = 10.84225142077871
$ws.Range("F6").Value = 25.56019909244887
$ws.Range("G6").Value = 34.81213511163889
$ws.Range("H6").Value = 13.526919237882
$ws.Range("I6").Value = 17.53365173822968
$ws.Range("M6").Value = 17.4494105548347
$ws.Range("N6").Value = 16.69850054216067

$ws.Range("B7").Value = 6.842510644523197
$ws.Range("D7").Value = 3.213155762846229
$ws.Range("E7").Value = 10.81293639081975
$ws.Range("F7").Value = 25.84525217626835
$ws.Range("G7").Value = 35.40555143632503
$ws.Range("H7").Value = 13.58704057692082
$ws.Range("I7").Value = 17.51116550462779
$ws.Range("M7").Value = 17.69291784332987
$ws.Range("N7").Value = 16.70403023771417

$ws.Range("B8").Value = 7.08546782846554
$ws.Range("D8").Value = 3.23179619608532
$ws.Range("E8").Value = 10.68970025017861
$ws.Range("F8").Value = 27.09072660808541
$ws.Range("G8").Value = 37.94242438528103
$ws.Range("H8").Value = 13.86065255094231
$ws.Range("I8").Value = 17.4203545683044
$ws.Range("M8").Value = 18.73645360723136
$ws.Range("N8").Value = 16.73314177615014

$ws.Range("B9").Value = 7.542001657786275
$ws.Range("D9").Value = 3.279069813196457
$ws.Range("E9").Value = 10.46983714849303
$ws.Range("F9").Value = 29.46092769062381
$ws.Range("G9").Value = 42.5788612755204
$ws.Range("H9").Value = 14.42250350450595
$ws.Range("I9").Value = 17.27251250553626
$ws.Range("M9").Value = 20.65290693831806
$ws.Range("N9").Value = 16.80680131513294

$ws.Range("B10").Value = 7.861959212412339
$ws.Range("D10").Value = 3.320263910469704
$ws.Range("E10").Value = 10.32139382752222
$ws.Range("F10").Value = 31.13699360560532
$ws.Range("G10").Value = 45.74849311264558
$ws.Range("H10").Value = 14.84653074450542
$ws.Range("I10").Value = 17.18231665668406
$ws.Range("M10").Value = 21.96844238940022
$ws.Range("N10").Value = 16.87068560099842

$ws.Range("B11").Value = 8.003545257199853
$ws.Range("D11").Value = 3.372193830603782
$ws.Range("E11").Value = 10.25666204977734
$ws.Range("F11").Value = 31.88151313408289
$ws.Range("G11").Value = 47.13427301494739
$ws.Range("H11").Value = 15.0409833283446
$ws.Range("I11").Value = 17.14528220591784
$ws.Range("M11").Value = 22.54463952185646
$ws.Range("N11").Value = 16.90179248740682

$ws.Range("B12").Value = 8.056544826091546
$ws.Range("D12").Value = 3.396577823217257
$ws.Range("E12").Value = 10.23254847280973
$ws.Range("F12").Value = 32.16059633441947
$ws.Range("G12").Value = 47.65065836262198
$ws.Range("H12").Value = 15.11476914077003
$ws.Range("I12").Value = 17.13183263332164
$ws.Range("M12").Value = 22.75948374120577
$ws.Range("N12").Value = 16.91385913076721

$ws.Range("B13").Value = 8.045158503519795
$ws.Range("D13").Value = 3.391342837851548
$ws.Range("E13").Value = 10.23772406933488
$ws.Range("F13").Value = 32.10062151887529
$ws.Range("G13").Value = 47.53982258135398
$ws.Range("H13").Value = 15.09887248764644
$ws.Range("I13").Value = 17.13470368666518
$ws.Range("M13").Value = 22.7133642951603
$ws.Range("N13").Value = 16.91124770660378

$ws.Range("B14").Value = 8.007918157103754
$ws.Range("D14").Value = 3.374207361762671
$ws.Range("E14").Value = 10.25467023365324
$ws.Range("F14").Value = 31.90453203996734
$ws.Range("G14").Value = 47.17692586911504
$ws.Range("H14").Value = 15.04705112296574
$ws.Range("I14").Value = 17.14416418674232
$ws.Range("M14").Value = 22.56238278489783
$ws.Range("N14").Value = 16.90277950366285

$ws.Range("B15").Value = 7.985025858246337
$ws.Range("D15").Value = 3.363663057622587
$ws.Range("E15").Value = 10.26510211066606
$ws.Range("F15").Value = 31.7840428847872
$ws.Range("G15").Value = 46.95354162360475
$ws.Range("H15").Value = 15.01532653401314
$ws.Range("I15").Value = 17.1500338396772
$ws.Range("M15").Value = 22.46946203316618
$ws.Range("N15").Value = 16.89762966704086

$ws.Range("B16").Value = 7.852622436054109
$ws.Range("D16").Value = 3.318976745513148
$ws.Range("E16").Value = 10.3256801963601
$ws.Range("F16").Value = 31.08795270780165
$ws.Range("G16").Value = 45.65677221663939
$ws.Range("H16").Value = 14.83384788195595
$ws.Range("I16").Value = 17.18481735219953
$ws.Range("M16").Value = 21.93032552938669
$ws.Range("N16").Value = 16.86869326107147

$ws.Range("B17").Value = 7.770347838971958
$ws.Range("D17").Value = 3.307849813375878
$ws.Range("E17").Value = 10.36355677013305
$ws.Range("F17").Value = 30.65612752861015
$ws.Range("G17").Value = 44.84663632358404
$ws.Range("H17").Value = 14.7228633843301
$ws.Range("I17").Value = 17.20717937252162
$ws.Range("M17").Value = 21.59377056136433
$ws.Range("N17").Value = 16.85146072209652

$ws.Range("B18").Value = 7.722655002614663
$ws.Range("D18").Value = 3.301579626600561
$ws.Range("E18").Value = 10.38560571628005
$ws.Range("F18").Value = 30.40607942049092
$ws.Range("G18").Value = 44.37540165879822
$ws.Range("H18").Value = 14.65917921957151
$ws.Range("I18").Value = 17.22041755251448
$ws.Range("M18").Value = 21.39810648159698
$ws.Range("N18").Value = 16.84174199772819

$ws.Range("B19").Value = 7.706444835626526
$ws.Range("D19").Value = 3.299479024110706
$ws.Range("E19").Value = 10.39311643406199
$ws.Range("F19").Value = 30.32113878713424
$ws.Range("G19").Value = 44.21495574550136
$ws.Range("H19").Value = 14.63764509198047
$ws.Range("I19").Value = 17.22496437602249
$ws.Range("M19").Value = 21.33150457580776
$ws.Range("N19").Value = 16.83848475853441

$ws.Range("B20").Value = 7.779144818927562
$ws.Range("D20").Value = 3.309020890388621
$ws.Range("E20").Value = 10.35949751162927
$ws.Range("F20").Value = 30.70227134811432
$ws.Range("G20").Value = 44.93342412073215
$ws.Range("H20").Value = 14.73466276544734
$ws.Range("I20").Value = 17.20475996884622
$ws.Range("M20").Value = 21.62981452475073
$ws.Range("N20").Value = 16.85327523308691

$ws.Range("B21").Value = 8.018873612519458
$ws.Range("D21").Value = 3.37925055031426
$ws.Range("E21").Value = 10.24968193193248
$ws.Range("F21").Value = 31.96220763104661
$ws.Range("G21").Value = 47.28374704962615
$ws.Range("H21").Value = 15.06226879682304
$ws.Range("I21").Value = 17.14136981556455
$ws.Range("M21").Value = 22.60682169128327
$ws.Range("N21").Value = 16.90525908422664

$ws.Range("B22").Value = 8.171941405100762
$ws.Range("D22").Value = 3.449528335817446
$ws.Range("E22").Value = 10.18023500140366
$ws.Range("F22").Value = 32.76893214937775
$ws.Range("G22").Value = 48.77087108740598
$ws.Range("H22").Value = 15.27722483257682
$ws.Range("I22").Value = 17.10328966379829
$ws.Range("M22").Value = 23.22577707259365
$ws.Range("N22").Value = 16.94090380318338

$ws.Range("B23").Value = 8.090590540036276
$ws.Range("D23").Value = 3.412219303498591
$ws.Range("E23").Value = 10.21708852308625
$ws.Range("F23").Value = 32.33997915171638
$ws.Range("G23").Value = 47.98173244746828
$ws.Range("H23").Value = 15.16244506857532
$ws.Range("I23").Value = 17.12330735648814
$ws.Range("M23").Value = 22.89726385845315
$ws.Range("N23").Value = 16.92172909038196

$ws.Range("B24").Value = 7.775168921983408
$ws.Range("D24").Value = 3.308491051160271
$ws.Range("E24").Value = 10.36133185121008
$ws.Range("F24").Value = 30.68141526953259
$ws.Range("G24").Value = 44.89420439222142
$ws.Range("H24").Value = 14.7293278787015
$ws.Range("I24").Value = 17.20585259153404
$ws.Range("M24").Value = 21.61352580925264
$ws.Range("N24").Value = 16.85245430502912

$ws.Range("B25").Value = 7.42098843883601
$ws.Range("D25").Value = 3.265138485433519
$ws.Range("E25").Value = 10.52700261313482
$ws.Range("F25").Value = 28.82979476913687
$ws.Range("G25").Value = 41.36429768869589
$ws.Range("H25").Value = 14.26823080901916
$ws.Range("I25").Value = 17.30927146607559
$ws.Range("M25").Value = 20.14985914725671
$ws.Range("N25").Value = 16.78513402192961
